$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace generic "Member N" placeholders with actual team member names
$ws.Range("A8").Value = "Lukas Hasler"
$ws.Range("A9").Value = "Pascal Strebel"
$ws.Range("A10").Value = "Cedric Weibel"
$ws.Range("A11").Value = "Robin Schmidiger"
$ws.Range("A12").Value = ""

# Update selection / view state
$ws.Range("D11").Select()
